{"js": "// Replace the date title and the 25 two-digit multiplication problems with\n// the new values from the commit. Every old value is unique in the document,\n// so a body-wide search-and-replace keyed on exact old text is unambiguous.\nconst replacements = [\n  [\"2024-03-06 Wednesday\", \"2024-03-07 Thursday\"],\n  [\"54\u00d732=\", \"32\u00d728=\"],\n  [\"33\u00d777=\", \"62\u00d775=\"],\n  [\"59\u00d780=\", \"25\u00d736=\"],\n  [\"48\u00d724=\", \"65\u00d786=\"],\n  [\"83\u00d764=\", \"62\u00d779=\"],\n  [\"70\u00d799=\", \"83\u00d731=\"],\n  [\"89\u00d744=\", \"69\u00d727=\"],\n  [\"27\u00d783=\", \"82\u00d774=\"],\n  [\"15\u00d771=\", \"45\u00d799=\"],\n  [\"16\u00d788=\", \"53\u00d759=\"],\n  [\"70\u00d772=\", \"45\u00d785=\"],\n  [\"81\u00d796=\", \"22\u00d798=\"],\n  [\"99\u00d778=\", \"84\u00d729=\"],\n  [\"81\u00d787=\", \"54\u00d792=\"],\n  [\"93\u00d755=\", \"96\u00d730=\"],\n  [\"79\u00d753=\", \"26\u00d733=\"],\n  [\"97\u00d793=\", \"32\u00d757=\"],\n  [\"28\u00d720=\", \"17\u00d738=\"],\n  [\"89\u00d740=\", \"73\u00d760=\"],\n  [\"68\u00d778=\", \"18\u00d783=\"],\n  [\"91\u00d760=\", \"63\u00d796=\"],\n  [\"69\u00d739=\", \"17\u00d715=\"],\n  [\"96\u00d765=\", \"42\u00d711=\"],\n  [\"27\u00d775=\", \"34\u00d788=\"],\n  [\"77\u00d748=\", \"18\u00d783=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date title and the 25 two-digit multiplication problems with\n# the new values from the commit. Every old value is unique in the document,\n# so Find/Replace (wdReplaceAll) keyed on the exact old text is unambiguous\n# and order-independent.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-03-06 Wednesday\", \"2024-03-07 Thursday\"),\n    @(\"54\u00d732=\", \"32\u00d728=\"),\n    @(\"33\u00d777=\", \"62\u00d775=\"),\n    @(\"59\u00d780=\", \"25\u00d736=\"),\n    @(\"48\u00d724=\", \"65\u00d786=\"),\n    @(\"83\u00d764=\", \"62\u00d779=\"),\n    @(\"70\u00d799=\", \"83\u00d731=\"),\n    @(\"89\u00d744=\", \"69\u00d727=\"),\n    @(\"27\u00d783=\", \"82\u00d774=\"),\n    @(\"15\u00d771=\", \"45\u00d799=\"),\n    @(\"16\u00d788=\", \"53\u00d759=\"),\n    @(\"70\u00d772=\", \"45\u00d785=\"),\n    @(\"81\u00d796=\", \"22\u00d798=\"),\n    @(\"99\u00d778=\", \"84\u00d729=\"),\n    @(\"81\u00d787=\", \"54\u00d792=\"),\n    @(\"93\u00d755=\", \"96\u00d730=\"),\n    @(\"79\u00d753=\", \"26\u00d733=\"),\n    @(\"97\u00d793=\", \"32\u00d757=\"),\n    @(\"28\u00d720=\", \"17\u00d738=\"),\n    @(\"89\u00d740=\", \"73\u00d760=\"),\n    @(\"68\u00d778=\", \"18\u00d783=\"),\n    @(\"91\u00d760=\", \"63\u00d796=\"),\n    @(\"69\u00d739=\", \"17\u00d715=\"),\n    @(\"96\u00d765=\", \"42\u00d711=\"),\n    @(\"27\u00d775=\", \"34\u00d788=\"),\n    @(\"77\u00d748=\", \"18\u00d783=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n\n$d.Save()\n"}
